$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 163, shifting existing rows 163:266 down to 164:267
$ws.Rows.Item(163).EntireRow.Insert()

# Populate the newly inserted row 163 with the new data record
$ws.Range("A163").Value = 10
$ws.Range("B163").Value = "Vega Modelo de Temuco"
$ws.Range("C163").Value = "La Araucanía"
$ws.Range("D163").Value = 44767
$ws.Range("E163").Value = 9
$ws.Range("F163").Value = 100112039
$ws.Range("G163").Value = "Ciboulette"
$ws.Range("H163").Value = "Sin especificar"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 50
$ws.Range("K163").Value = 7000
$ws.Range("L163").Value = 8000
$ws.Range("M163").Value = 7600
$ws.Range("N163").Value = "$/docena de atados"
$ws.Range("O163").Value = "Provincia de Cautín"
$ws.Range("P163").Value = 2533
$ws.Range("Q163").Value = 3
$ws.Range("R163").Value = "Hortaliza"
